$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 5.907218141265402
$ws.Range("E4").Value = 5.917486466529609

$ws.Range("C5").Value = 6.511263427347003
$ws.Range("E5").Value = 6.422943767670297

$ws.Range("C6").Value = 4.325828829470257
$ws.Range("E6").Value = 4.493586323244281

$ws.Range("C7").Value = 2.964652118442834
$ws.Range("E7").Value = 3.83627393798931

$ws.Range("C8").Value = 3.171852776411788
$ws.Range("E8").Value = 4.030605385534614

$ws.Range("C9").Value = 1.773712379859993
$ws.Range("E9").Value = 3.484530515673856

$ws.Range("C10").Value = 2.533350906619081
$ws.Range("E10").Value = 3.524103740130435

$ws.Range("C11").Value = 2.661040979345697
$ws.Range("E11").Value = 3.567108445582057

$ws.Range("C12").Value = 3.150198973767537
$ws.Range("E12").Value = 3.699072253610103

$ws.Range("C13").Value = 0.4641929091049102
$ws.Range("E13").Value = 2.550259844884462

$ws.Range("C14").Value = 2.585454129751663
$ws.Range("E14").Value = 2.671828487424377

$ws.Range("C15").Value = -0.4532848472497908
$ws.Range("E15").Value = 2.066462658785673

$ws.Range("C16").Value = 0.9477102747197819
$ws.Range("E16").Value = 1.83067479293082

$ws.Range("C17").Value = 2.222852754198135
$ws.Range("E17").Value = 1.934107558751452

$ws.Range("C18").Value = -0.007094633234694392
$ws.Range("E18").Value = 1.444584248586422

$ws.Range("C19").Value = 3.078872076370009
$ws.Range("E19").Value = 2.279508996785351
